$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "customers"
$ws.Range("C1").Value = "Sales"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "dewf"
$ws.Range("C2").Value = 32423

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "fewrf"
$ws.Range("C3").Value = 43255

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ve"
$ws.Range("C4").Value = 5235

# Match the selection recorded in the diff (active cell C4)
$ws.Range("C4").Select()
